# Adds the "1st March" meeting-minutes section to the end of the document,
# reproducing the bullet list (with a nested sub-list and a spell-check
# exception run around "Transfuzer") exactly as it appears in the source
# revision. We build the new paragraphs as a literal WordprocessingML
# fragment and drop it in right after the existing last paragraph (the one
# holding the "_GoBack" bookmark), which keeps that paragraph, its bookmark,
# and all prior content completely untouched.

$d = $word.ActiveDocument

# Collapsed range positioned at the very end of the document's content,
# i.e. immediately after the final paragraph mark of the last paragraph
# ("We need to meet up before Wednesday to practice presentation").
$endPos = $d.Content.End
$insertionRange = $d.Range($endPos, $endPos)

$newSectionXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>st</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> March</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="none"/></w:rPr><w:t>Presentation day</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="none"/></w:rPr><w:t>Lee was unable to make it to the presentation due to flu</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve">We need to have a plan for if we get onto the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="none"/></w:rPr><w:t>Transfuzer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve"> Project. Write down what we plan to do and such</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="none"/></w:rPr><w:t>We need to get people playtesting the game, even in its unfinished state we need to see how players react to the environment</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="none"/></w:rPr><w:t>Set up next sprint as a group (Lee was still unable to make it in).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$countBefore = $d.Paragraphs.Count
$insertionRange.InsertXML($newSectionXml)
$countAfter = $d.Paragraphs.Count

Write-Output ("Paragraphs before: " + $countBefore + ", after: " + $countAfter)
Write-Output ("Last paragraph text: " + $d.Paragraphs.Last.Range.Text)
